$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: B6 -- new "esriVersion" code-style parameter line ---
$b6 = $ws.Range("B6")
$b6.Value = 'let esriVersion = "4.28";     '

# "let " (chars 1-4) is left with the cell's base/default formatting (no
# explicit run formatting), matching the other code-style rows on the sheet.

# "esriVersion" (chars 5-15): bold, purple
$r1 = $b6.Characters(5, 11)
$r1.Font.Name = "Arial Unicode MS"
$r1.Font.Size = 10
$r1.Font.Bold = $true
$r1.Font.Color = 0x910083

# " " (char 16): purple, not bold
$r2 = $b6.Characters(16, 1)
$r2.Font.Name = "Arial Unicode MS"
$r2.Font.Size = 10
$r2.Font.Color = 0x910083

# "= " (chars 17-18): dark gray/near-black
$r3 = $b6.Characters(17, 2)
$r3.Font.Name = "Arial Unicode MS"
$r3.Font.Size = 10
$r3.Font.Color = 0x080808

# "4.28" quoted (chars 19-24): green (string literal)
$r4 = $b6.Characters(19, 6)
$r4.Font.Name = "Arial Unicode MS"
$r4.Font.Size = 10
$r4.Font.Color = 0x177D06

# ";     " (chars 25-30): dark gray/near-black
$r5 = $b6.Characters(25, 6)
$r5.Font.Name = "Arial Unicode MS"
$r5.Font.Size = 10
$r5.Font.Color = 0x080808

# Give B6 the same vertical-centered (no forced horizontal) alignment used
# for the other parameter-name cells in column B.
$b6.VerticalAlignment = -4108

# --- Row 6: C6 -- description of the new parameter ---
$c6 = $ws.Range("C6")
$c6.Value = "Version of the ESRI SDK -- ESRI puts out a new version 3-4 times each year"

# Italicize the explanatory clause at the end of the description.
$desc = $c6.Characters(28, 49)
$desc.Font.Name = "Aptos Narrow"
$desc.Font.Size = 11
$desc.Font.Italic = $true

# Move the active selection down to C7, below the newly added row.
$ws.Range("C7").Select()
